# Apply "New crime data collected" weekly update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
# Volume/Number text: 46 -> 47
$ws.Range("C1").Value = "Volume 30   Number  47"

# Report covering the week dates
$ws.Range("C6").Value = "Report Covering the Week  11/20/2023  Through  11/26/2023"

# --- Row 14 (Murder) ---
$ws.Range("N14").Value = -95.238095238095

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = "0"
$ws.Range("H15").Value = "***.*"
$ws.Range("N15").Value = -30

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = -45.833333333333
$ws.Range("I16").Value = 192
$ws.Range("J16").Value = 206
$ws.Range("K16").Value = -6.796116504854
$ws.Range("L16").Value = 10.344827586206
$ws.Range("M16").Value = -20.331950207468
$ws.Range("N16").Value = -75.909661229611

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 42.105263157894
$ws.Range("I17").Value = 328
$ws.Range("J17").Value = 320
$ws.Range("K17").Value = 2.5
$ws.Range("L17").Value = 3.470031545741
$ws.Range("M17").Value = 144.776119402985
$ws.Range("N17").Value = -13.684210526315

# --- Row 18 (Burglary) ---
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 125
$ws.Range("J18").Value = 139
$ws.Range("K18").Value = -10.071942446043
$ws.Range("L18").Value = 25
$ws.Range("M18").Value = -54.212454212454
$ws.Range("N18").Value = -88.594890510948

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -21.568627450980
$ws.Range("I19").Value = 562
$ws.Range("J19").Value = 605
$ws.Range("K19").Value = -7.107438016528
$ws.Range("L19").Value = 21.120689655172
$ws.Range("M19").Value = 71.865443425076
$ws.Range("N19").Value = 5.243445692883

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 22.222222222222
$ws.Range("I20").Value = 236
$ws.Range("J20").Value = 202
$ws.Range("K20").Value = 16.831683168316
$ws.Range("L20").Value = 24.867724867724
$ws.Range("M20").Value = -14.801444043321
$ws.Range("N20").Value = -92.107023411371

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -13.636363636363
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 116
$ws.Range("H21").Value = -5.172413793103
$ws.Range("I21").Value = 1465
$ws.Range("J21").Value = 1495
$ws.Range("K21").Value = -2.006688963210
$ws.Range("L21").Value = 16.177636796193
$ws.Range("M21").Value = 14.992150706436
$ws.Range("N21").Value = -74.948700410396

# --- Row 22 (Transit) ---
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("M22").Value = -30.769230769230

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = 42.307692307692
$ws.Range("I24").Value = 1203
$ws.Range("J24").Value = 1227
$ws.Range("K24").Value = -1.955990220048
$ws.Range("L24").Value = 37.172177879133
$ws.Range("M24").Value = 91.255961844197

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 62.5
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 489
$ws.Range("J25").Value = 484
$ws.Range("K25").Value = 1.033057851239
$ws.Range("L25").Value = 15.058823529411
$ws.Range("M25").Value = 7.947019867549

# --- Row 26 (UCR Rape*) ---
$ws.Range("D26").Value = "0"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 34
$ws.Range("K26").Value = -2.857142857142
$ws.Range("L26").Value = 88.888888888888

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = -25.925925925925
$ws.Range("L27").Value = 8.108108108108
